$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
$ws.Range("A8").Value = "Volume 32   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/19/2025  Through  5/25/2025"

# --- Crime statistics table updates (rows 16-31) ---
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 5
$ws.Range("H16").Value = -28.571428571428
$ws.Range("I16").Value = 32
$ws.Range("J16").Value = 32
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = -21.951219512195
$ws.Range("N16").Value = -88.405797101449
$ws.Range("C17").Value = 2
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("D17").Value = 5
$ws.Range("D17").NumberFormat = "#,##0"
$ws.Range("E17").Value = -60
$ws.Range("E17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 28.571428571428
$ws.Range("I17").Value = 63
$ws.Range("J17").Value = 41
$ws.Range("K17").Value = 53.658536585365
$ws.Range("L17").Value = 50
$ws.Range("M17").Value = 231.578947368421
$ws.Range("N17").Value = 26
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 166.666666666667
$ws.Range("I18").Value = 65
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = 62.5
$ws.Range("L18").Value = 22.641509433962
$ws.Range("M18").Value = 38.297872340425
$ws.Range("N18").Value = -87.758945386064
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 10
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = -39.534883720930
$ws.Range("I19").Value = 148
$ws.Range("J19").Value = 174
$ws.Range("K19").Value = -14.942528735632
$ws.Range("L19").Value = -21.276595744680
$ws.Range("M19").Value = 3.496503496503
$ws.Range("N19").Value = -61.154855643044
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 30
$ws.Range("I20").Value = 64
$ws.Range("J20").Value = 56
$ws.Range("K20").Value = 14.285714285714
$ws.Range("L20").Value = 42.222222222222
$ws.Range("M20").Value = 33.333333333333
$ws.Range("N20").Value = -95.395683453237
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -23.809523809523
$ws.Range("F21").Value = 61
$ws.Range("G21").Value = 71
$ws.Range("H21").Value = -14.084507042253
$ws.Range("I21").Value = 378
$ws.Range("J21").Value = 347
$ws.Range("K21").Value = 8.933717579250
$ws.Range("L21").Value = 4.419889502762
$ws.Range("M21").Value = 26.421404682274
$ws.Range("N21").Value = -85.638297872340
$ws.Range("D22").Value = 2
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -83.333333333333
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = -17.647058823529
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = 2.857142857142
$ws.Range("F24").Value = 152
$ws.Range("G24").Value = 112
$ws.Range("H24").Value = 35.714285714285
$ws.Range("I24").Value = 806
$ws.Range("J24").Value = 635
$ws.Range("K24").Value = 26.929133858267
$ws.Range("L24").Value = 26.929133858267
$ws.Range("M24").Value = 109.350649350649
$ws.Range("C25").Value = 30
$ws.Range("D25").Value = 22
$ws.Range("E25").Value = 36.363636363636
$ws.Range("F25").Value = 111
$ws.Range("G25").Value = 70
$ws.Range("H25").Value = 58.571428571428
$ws.Range("I25").Value = 633
$ws.Range("J25").Value = 470
$ws.Range("K25").Value = 34.680851063829
$ws.Range("L25").Value = 40.979955456570
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 83.333333333333
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = 15
$ws.Range("I26").Value = 128
$ws.Range("J26").Value = 107
$ws.Range("K26").Value = 19.626168224299
$ws.Range("L26").Value = 31.958762886597
$ws.Range("M26").Value = 43.820224719101
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 12
$ws.Range("J28").Value = 10
$ws.Range("K28").Value = 20
$ws.Range("L28").Value = -7.692307692307
$ws.Range("L31").Value = -36.363636363636
